# Update countries & provincias Spain
# Applies the 8-Sep-2020 11:31 data refresh to the "Pais" sheet:
#  - bumps the "last updated" timestamp in A1
#  - refreshes case numbers for several countries
#  - two country pairs overtake their neighbour in Casos totales and swap
#    places in the (already B-descending-sorted) table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- timestamp header -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 8 de Septiembre de 2020 a las 11:31"

# --- straightforward data refresh (no reordering) ----------------------
# Indonesia (row 26)
Set-Row 26 200035 3046 142958 48847 0 100 8230

# Israel (row 28)
Set-Row 28 135288 1313 106297 27960 0 5 1031

# Polonia (row 49) - only Casos activos / Recuperados change
$ws.Cells.Item(49, 4).Value = 55910
$ws.Cells.Item(49, 5).Value = 13092

# Croacia (row 90)
Set-Row 90 12285 204 9553 2529 0 2 203

# --- Austria overtakes Irlanda (rows 70/71) -----------------------------
$ws.Range("A70").Value = "Austria"
Set-Row 70 30081 520 25629 3705 0 1 747

$ws.Range("A71").Value = "Irlanda"
Set-Row 71 29774 0 23364 4633 0 0 1777

# --- Hong Kong overtakes Congo (rows 113/114) ---------------------------
$ws.Range("A113").Value = "Hong Kong"
Set-Row 113 4896 6 4543 254 0 1 99

$ws.Range("A114").Value = "Congo"
Set-Row 114 4891 0 3887 902 0 0 102

# --- Eslovaquia overtakes Nicaragua (rows 117/118) ----------------------
$ws.Range("A117").Value = "Eslovaquia"
Set-Row 117 4727 91 2913 1777 0 0 37

$ws.Range("A118").Value = "Nicaragua"
Set-Row 118 4668 0 2913 1614 0 0 141

# --- Birmania overtakes Guyana (rows 153/154) ---------------------------
$ws.Range("A153").Value = "Birmania"
Set-Row 153 1610 92 388 1214 0 0 8

$ws.Range("A154").Value = "Guyana"
Set-Row 154 1560 0 962 551 0 0 47
